# Sample Project / Main.xlsx - "Rules" sheet edit
#
# Target change: cell B11 changes from the shared string "R40" to the
# text value "1" (still a text/string cell, not a number).
#
# A bare "1" typed (or assigned) into a General-formatted cell is
# auto-detected by Excel as a number, which would store the cell as
# numeric instead of text. Prefixing it with an apostrophe is the normal
# Excel way to force text entry, but it also marks the cell with the
# "number stored as text" (quote-prefix) formatting. To keep B11's
# original formatting untouched, its current format is first copied to
# an out-of-the-way scratch cell, then pasted back onto B11 once the new
# text value has been entered; finally the scratch cell is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")
$scratch = $ws.Range("ZZ500")

# Remember B11's current formatting.
$target.Copy($scratch)

# Enter the new value, forced to text so it stays a string (not a number).
$target.Value = "'1"

# Re-apply B11's original formatting, then tidy up the scratch cell.
$scratch.Copy()
$target.PasteSpecial(-4122)
$scratch.Clear()
